$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.656.18'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.042.26'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '542.19'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '133.37'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -10.30%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.036.74'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.99%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.486'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.36'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -10.14%  '
$ws.Range('E11').Value = '  -3.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.455'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '34.50'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000213'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.63%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.526.16'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '62.687.45'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.039.93'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.55'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '478.35'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -10.59%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.28'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.49%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.691'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.96'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -6.40%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '76.82'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.64%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.09'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.68'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.15'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -7.12%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.91'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -10.76%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.12'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '60.26'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +13.67%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.46'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.90%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '507.97'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -8.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.87'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.52%  '
$ws.Range('E37').Value = '  -6.83%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0395'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -11.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.033.42'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0782'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.117'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.97'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.53'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -11.72%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.249'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.00'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.00%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '118.60'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '24.21'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.106'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.00%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₃0487'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -7.25%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.31'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +57.96%  '
